$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 686
$ws.Cells.Item(686, 1).Value = 'Emon'
$ws.Cells.Item(686, 2).Value = 'conccn'
$ws.Cells.Item(686, 3).Value = 1
$ws.Cells.Item(686, 4).Value = 'longitude latitude alevel time'
$ws.Cells.Item(686, 5).Value = 'Aerosol Number Concentration'
$ws.Cells.Item(686, 6).Value = 'm-3'
$ws.Cells.Item(686, 7).Value = 'web'
$ws.Cells.Item(686, 8).Value = 'Available in TM5, though yet to be added by Tommi'
$ws.Cells.Item(686, 9).Value = 'Tommi Bergman, Thomas'
$ws.Cells.Item(686, 10).Value = '''''Number concentration'' means the number of particles or other specified objects per unit volume. ''Aerosol'' means the system of suspended liquid or solid particles in air (except cloud droplets) and their carrier gas, the air itself. ''Ambient_aerosol'' means that the aerosol is measured or modelled at the ambient state of pressure, temperature and relative humidity that exists in its immediate environment. ''Ambient aerosol particles'' are aerosol particles that have taken up ambient water through hygroscopic growth. The extent of hygroscopic growth depends on the relative humidity and the composition of the particles.'
$ws.Cells.Item(686, 11).Value = 'GeoMIP,VIACSAB'

# Row 687
$ws.Cells.Item(687, 1).Value = 'Emon'
$ws.Cells.Item(687, 2).Value = 'sconcss'
$ws.Cells.Item(687, 3).Value = 1
$ws.Cells.Item(687, 4).Value = 'longitude latitude time'
$ws.Cells.Item(687, 5).Value = 'Surface Concentration of Seasalt'
$ws.Cells.Item(687, 6).Value = 'kg m-3'
$ws.Cells.Item(687, 7).Value = 'web'
$ws.Cells.Item(687, 8).Value = 'Available in TM5, though yet to be added by Tommi'
$ws.Cells.Item(687, 9).Value = 'Tommi Bergman, Thomas'
$ws.Cells.Item(687, 10).Value = 'mass concentration of seasalt dry aerosol in air in model lowest layer'
$ws.Cells.Item(687, 11).Value = 'GeoMIP,VIACSAB'

# Row 691
$ws.Cells.Item(691, 1).Value = 'Omon'
$ws.Cells.Item(691, 2).Value = 'talkos'
$ws.Cells.Item(691, 3).Value = 1
$ws.Cells.Item(691, 4).Value = 'longitude latitude time'
$ws.Cells.Item(691, 5).Value = 'Surface Total Alkalinity'
$ws.Cells.Item(691, 6).Value = 'mol m-3'
$ws.Cells.Item(691, 7).Value = 'web'
$ws.Cells.Item(691, 8).Value = 'AlkaliniSFC_E3T / e3t'
$ws.Cells.Item(691, 9).Value = 'Raffaele Bernardello'
$ws.Cells.Item(691, 10).Value = 'total alkalinity equivalent concentration (including carbonate, borate, phosphorus, silicon, and nitrogen components)'
$ws.Cells.Item(691, 11).Value = 'AerChemMIP,C4MIP,CMIP,GMMIP,GeoMIP,HighResMIP,LS3MIP,OMIP'

# Row 692
$ws.Cells.Item(692, 1).Value = 'Omon'
$ws.Cells.Item(692, 2).Value = 'phos'
$ws.Cells.Item(692, 3).Value = 2
$ws.Cells.Item(692, 4).Value = 'longitude latitude time'
$ws.Cells.Item(692, 5).Value = 'Surface pH'
$ws.Cells.Item(692, 6).Value = 1
$ws.Cells.Item(692, 7).Value = 'web'
$ws.Cells.Item(692, 8).Value = 'PHSFC'
$ws.Cells.Item(692, 9).Value = 'Raffaele Bernardello'
$ws.Cells.Item(692, 10).Value = 'negative log10 of hydrogen ion concentration with the concentration expressed as mol H kg-1.'
$ws.Cells.Item(692, 11).Value = 'AerChemMIP,C4MIP,CMIP,GMMIP,GeoMIP,HighResMIP,LS3MIP,OMIP'

# Row 693
$ws.Cells.Item(693, 1).Value = 'Omon'
$ws.Cells.Item(693, 2).Value = 'po4os'
$ws.Cells.Item(693, 3).Value = 1
$ws.Cells.Item(693, 4).Value = 'longitude latitude time'
$ws.Cells.Item(693, 5).Value = 'mole_concentration_of_dissolved_inorganic_phosphorous_in_sea_water'
$ws.Cells.Item(693, 6).Value = 'mol m-3'
$ws.Cells.Item(693, 7).Value = 'web'
$ws.Cells.Item(693, 8).Value = 'PO4_E3T/e3t'
$ws.Cells.Item(693, 9).Value = 'Raffaele Bernardello'
$ws.Cells.Item(693, 10).Value = 'Mole concentration means number of moles per unit volume, also called ''molarity'', and is used in the construction ''mole_concentration_of_X_in_Y'', where X is a material constituent of Y. A chemical or biological species denoted by X may be described by a single term such as ''nitrogen'' or a phrase such as ''nox_expressed_as_nitrogen''. ''Dissolved inorganic phosphorus'' means the sum of all inorganic phosphorus in solution (including phosphate, hydrogen phosphate, dihydrogen phosphate, and phosphoric acid).'
$ws.Cells.Item(693, 11).Value = 'AerChemMIP,C4MIP,CMIP,GMMIP,GeoMIP,HighResMIP,LS3MIP,OMIP'

# Row 694
$ws.Cells.Item(694, 1).Value = 'Omon'
$ws.Cells.Item(694, 2).Value = 'wfo'
$ws.Cells.Item(694, 3).Value = 1
$ws.Cells.Item(694, 4).Value = 'longitude latitude time'
$ws.Cells.Item(694, 5).Value = 'Water Flux into Sea Water'
$ws.Cells.Item(694, 6).Value = 'kg m-2 s-1'
$ws.Cells.Item(694, 7).Value = 'web'
$ws.Cells.Item(694, 8).Value = 'wfo'
$ws.Cells.Item(694, 9).Value = 'Raffaele Bernardello'
$ws.Cells.Item(694, 10).Value = 'computed as the water  flux into the ocean divided by the area of the ocean portion of the grid cell.  This is the sum of the next two variables in this table.'
$ws.Cells.Item(694, 11).Value = 'AerChemMIP,C4MIP,CMIP,DAMIP,GMMIP,GeoMIP,HighResMIP,LS3MIP,OMIP,VIACSAB,VolMIP'

# Row 695
$ws.Cells.Item(695, 1).Value = 'Omon'
$ws.Cells.Item(695, 2).Value = 'zhalfo'
$ws.Cells.Item(695, 3).Value = 1
$ws.Cells.Item(695, 4).Value = 'longitude latitude olevel time'
$ws.Cells.Item(695, 5).Value = 'Depth Below Geoid of Interfaces Between Ocean Layers'
$ws.Cells.Item(695, 6).Value = 'm'
$ws.Cells.Item(695, 7).Value = 'web'
$ws.Cells.Item(695, 8).Value = 'tpt_dep'
$ws.Cells.Item(695, 9).Value = 'Raffaele Bernardello'
$ws.Cells.Item(695, 10).Value = 'Depth below geoid'
$ws.Cells.Item(695, 11).Value = 'AerChemMIP,C4MIP,CFMIP,CMIP,DAMIP,GMMIP,GeoMIP,HighResMIP,LS3MIP,OMIP,VIACSAB'

# Row 696
$ws.Cells.Item(696, 1).Value = 'Omon'
$ws.Cells.Item(696, 2).Value = 'intppcalc'
$ws.Cells.Item(696, 3).Value = 3
$ws.Cells.Item(696, 4).Value = 'longitude latitude time'
$ws.Cells.Item(696, 5).Value = 'Net Primary Mole Productivity of Carbon by Calcareous Phytoplankton'
$ws.Cells.Item(696, 6).Value = 'mol m-2 s-1'
$ws.Cells.Item(696, 7).Value = 'web'
$ws.Cells.Item(696, 8).Value = 'INTPCAL'
$ws.Cells.Item(696, 9).Value = 'Raffaele Bernardello'
$ws.Cells.Item(696, 10).Value = '''''Production of carbon'' means the production of biomass expressed as the mass of carbon which it contains. Net primary production is the excess of gross primary production (rate of synthesis of biomass from inorganic precursors) by autotrophs (''producers''), for example, photosynthesis in plants or phytoplankton, over the rate at which the autotrophs themselves respire some of this biomass. ''Productivity'' means production per unit area. Phytoplankton are autotrophic prokaryotic or eukaryotic algae that live near the water surface where there is sufficient light to support photosynthesis. ''Calcareous phytoplankton'' are phytoplankton that produce calcite. The phrase ''expressed_as'' is used in the construction A_expressed_as_B, where B is a chemical constituent of A. It means that the quantity indicated by the standard name is calculated solely with respect to the B contained in A, neglecting all other chemical constituents of A. Calcite is a mineral that is a polymorph of calcium carbonate.'
$ws.Cells.Item(696, 11).Value = 'AerChemMIP,C4MIP,CMIP,GMMIP,GeoMIP,HighResMIP,LS3MIP,OMIP,VIACSAB'

# Row 697
$ws.Cells.Item(697, 1).Value = 'Omon'
$ws.Cells.Item(697, 2).Value = 'intpcalcite'
$ws.Cells.Item(697, 3).Value = 3
$ws.Cells.Item(697, 4).Value = 'longitude latitude time'
$ws.Cells.Item(697, 5).Value = 'Calcite Production'
$ws.Cells.Item(697, 6).Value = 'mol m-2 s-1'
$ws.Cells.Item(697, 7).Value = 'web'
$ws.Cells.Item(697, 8).Value = 'INTPCAL'
$ws.Cells.Item(697, 9).Value = 'Raffaele Bernardello'
$ws.Cells.Item(697, 10).Value = 'Vertically integrated calcite production'
$ws.Cells.Item(697, 11).Value = 'AerChemMIP,C4MIP,CMIP,GMMIP,GeoMIP,HighResMIP,LS3MIP,OMIP,VIACSAB'

# Row 699
$ws.Cells.Item(699, 1).Value = 'Oclim'
$ws.Cells.Item(699, 2).Value = 'zhalfo'
$ws.Cells.Item(699, 3).Value = 1
$ws.Cells.Item(699, 4).Value = 'longitude latitude olevel time2'
$ws.Cells.Item(699, 5).Value = 'Depth Below Geoid of Interfaces Between Ocean Layers'
$ws.Cells.Item(699, 6).Value = 'm'
$ws.Cells.Item(699, 7).Value = 'web'
$ws.Cells.Item(699, 8).Value = 'tpt_dep'
$ws.Cells.Item(699, 9).Value = 'Raffaele Bernardello, Thomas'
$ws.Cells.Item(699, 10).Value = 'Depth below geoid'
$ws.Cells.Item(699, 11).Value = 'CMIP,FAFMIP,HighResMIP,LUMIP,RFMIP'

# Row 702
$ws.Cells.Item(702, 1).Value = 'Omon'
$ws.Cells.Item(702, 2).Value = 'spco2nat'
$ws.Cells.Item(702, 3).Value = 1
$ws.Cells.Item(702, 4).Value = 'longitude latitude time depth0m'
$ws.Cells.Item(702, 5).Value = 'Natural Surface Aqueous Partial Pressure of CO2'
$ws.Cells.Item(702, 6).Value = 'Pa'
$ws.Cells.Item(702, 7).Value = 'web'
$ws.Cells.Item(702, 8).Value = 'pCO2sea (in uatm) in simulation where ocean biogeochemistry sees preindustrial atmospheric pCO2 but initial conditions and forcings  are identical to historical. This variable will be delivered for OMIP but seems unlikely to be delivered for C4MIP. This is because it requires running twice simulations.'
$ws.Cells.Item(702, 9).Value = 'Raffaele Bernardello'
$ws.Cells.Item(702, 10).Value = 'The surface called ''surface'' means the lower boundary of the atmosphere. The chemical formula for carbon dioxide is CO2. In ocean biogeochemistry models, a ''natural analogue'' is used to simulate the effect on a modelled variable of imposing preindustrial atmospheric carbon dioxide concentrations, even when the model as a whole may be subjected to varying forcings. The partial pressure of a gaseous constituent of air is the pressure which it alone would exert with unchanged temperature and number of moles per unit volume. The partial pressure of a dissolved gas in sea water is the partial pressure in air with which it would be in equilibrium. The partial pressure difference between sea water and air is positive when the partial pressure of the dissolved gas in sea water is greater than the partial pressure in air.'
$ws.Cells.Item(702, 11).Value = 'AerChemMIP,C4MIP,CMIP,GMMIP,GeoMIP,HighResMIP,LS3MIP,OMIP'

# Row 704
$ws.Cells.Item(704, 1).Value = 'Oyr'
$ws.Cells.Item(704, 2).Value = 'dissicnat'
$ws.Cells.Item(704, 3).Value = 1
$ws.Cells.Item(704, 4).Value = 'longitude latitude olevel time'
$ws.Cells.Item(704, 5).Value = 'Natural Dissolved Inorganic Carbon Concentration'
$ws.Cells.Item(704, 6).Value = 'mol m-3'
$ws.Cells.Item(704, 7).Value = 'web'
$ws.Cells.Item(704, 8).Value = 'DIC_E3T/e3t in simulation where ocean biogeochemistry sees preindustrial atmospheric pCO2 but initial conditions and forcings  are identical to historical. This variable will be delivered for OMIP but seems unlikely to be delivered for C4MIP. This is because it requires running twice simulations.'
$ws.Cells.Item(704, 9).Value = 'Raffaele Bernardello'
$ws.Cells.Item(704, 10).Value = 'Dissolved inorganic carbon (CO3+HCO3+H2CO3) concentration at preindustrial atmospheric xCO2'
$ws.Cells.Item(704, 11).Value = 'AerChemMIP,CMIP,GeoMIP,LUMIP,OMIP'

# Row 705
$ws.Cells.Item(705, 1).Value = 'Oyr'
$ws.Cells.Item(705, 2).Value = 'phnat'
$ws.Cells.Item(705, 3).Value = 1
$ws.Cells.Item(705, 4).Value = 'longitude latitude olevel time'
$ws.Cells.Item(705, 5).Value = 'Natural pH'
$ws.Cells.Item(705, 6).Value = 1
$ws.Cells.Item(705, 7).Value = 'web'
$ws.Cells.Item(705, 8).Value = 'PH in simulation where ocean biogeochemistry sees preindustrial atmospheric pCO2 but initial conditions and forcings  are identical to historical. This variable will be delivered for OMIP but seems unlikely to be delivered for C4MIP. This is because it requires running twice simulations.'
$ws.Cells.Item(705, 9).Value = 'Raffaele Bernardello'
$ws.Cells.Item(705, 10).Value = 'negative log10 of hydrogen ion concentration with the concentration expressed as mol H kg-1.'
$ws.Cells.Item(705, 11).Value = 'AerChemMIP,CMIP,GeoMIP,LUMIP,OMIP'

$ws.Range("A684").Select()
